$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Metadata")

# Update Version value (row 3)
$ws1.Range("B3").Value = "6.0.0"

# Update Date value (row 8)
$ws1.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Update Publisher value (row 9)
$ws1.Range("B9").Value = "Alvearie Team"

# Replace the duplicated "Contact" / "No display for ContactDetail" rows (10-11)
# with a single "Jurisdiction" / "United States of America" row, then shift
# everything below up by one row (collapsing the duplicate row out).
$ws1.Range("A10").Value = "Jurisdiction"
$ws1.Range("B10").Value = "United States of America"

$ws1.Range("A11").Value = "Description"
$ws1.Range("B11").Value = "Supporting information type codes for claims"

$ws1.Range("A12").Value = "Purpose"
$ws1.Range("B12").Value = ""

$ws1.Range("A13").Value = "Copyright"
$ws1.Range("B13").Value = ""

$ws1.Range("A14").Value = "Immutable"
$ws1.Range("B14").Value = "BooleanType[null]"

# Delete now-duplicate last row (15)
$ws1.Rows.Item(15).Delete()
